$d = $word.ActiveDocument

# The first two paragraphs ("Day After Day - May 1937" heading, wrapped in a
# bookmark, and "By Dorothy Day" in bold) are replaced by a pandoc-style
# title block: a Title-styled paragraph (word-by-word runs, no bookmark) and
# an Authors-styled paragraph (word-by-word runs, no bold, no "By " prefix).

$titleWords = @("Day", " ", "After", " ", "Day", " ", "-", " ", "May", " ", "1937")
$authorWords = @("Dorothy", " ", "Day")

$titleRuns = ($titleWords | ForEach-Object { "<w:r><w:t xml:space=`"preserve`">$_</w:t></w:r>" }) -join ""
$authorRuns = ($authorWords | ForEach-Object { "<w:r><w:t xml:space=`"preserve`">$_</w:t></w:r>" }) -join ""

$bodyXml = "<w:p><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$titleRuns</w:p>" + `
           "<w:p><w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>$authorRuns</w:p>"

$packageXml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" + `
  "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" + `
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" + `
      "<pkg:xmlData>" + `
        "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" + `
          "<w:body>$bodyXml</w:body>" + `
        "</w:document>" + `
      "</pkg:xmlData>" + `
    "</pkg:part>" + `
  "</pkg:package>"

# The first paragraph (title, currently wrapped by the bookmark) and the
# second paragraph ("By Dorothy Day") together form the span to replace.
$rngStart = $d.Paragraphs(1).Range.Start
$rngEnd = $d.Paragraphs(2).Range.End
$rng = $d.Range($rngStart, $rngEnd)
$rng.InsertXML($packageXml)
